# ---------------------------------------------------------------------------
# Adds a new "2022-Q4" quarter sheet (duplicating the structure of the
# existing per-quarter sheets) and inserts its summary row at the top of the
# "总计" (total) sheet, pushing the existing rows down by one.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item(1)          # "总计"

# ---------------------------------------------------------------------------
# 1. "总计" sheet: insert a new row for 2022-Q4 at the top of the data
#    (row 2), shifting the existing data rows down by one.
# ---------------------------------------------------------------------------

# Extend the style of the last existing data row's index cell (A6, which
# carries the shared "index" style) down onto the new last row (A7) before
# shifting values, so every row in the resulting A2:A7 range keeps the same
# look as before.
$total.Cells.Item(6, 1).Copy($total.Cells.Item(7, 1))

# Shift columns B (quarter label), C (count) and D (market value) down one
# row at a time, working from the bottom up so we never overwrite data we
# still need to read.
for ($r = 6; $r -ge 2; $r--) {
    $newRow = $r + 1
    $total.Cells.Item($newRow, 2).Value = $total.Cells.Item($r, 2).Value2
    $total.Cells.Item($newRow, 3).Value = $total.Cells.Item($r, 3).Value2
    $total.Cells.Item($newRow, 4).Value = $total.Cells.Item($r, 4).Value2
}

# Write the new 2022-Q4 summary row.
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 25
$total.Cells.Item(2, 4).Value = 4.03

# Column A is a plain sequential row index (0-based) - recompute it for
# every data row now that a row has been added.
for ($r = 2; $r -le 7; $r++) {
    $total.Cells.Item($r, 1).Value = $r - 2
}

# ---------------------------------------------------------------------------
# 2. Add the new "2022-Q4" worksheet right after "总计", mirroring the
#    layout of the other quarterly sheets.
# ---------------------------------------------------------------------------

$newSheet = $wb.Worksheets.Add($null, $total)
$newSheet.Name = "2022-Q4"

# NOTE: sheet references in this host resolve by *position*, not identity -
# grab "2022-Q3" (used purely as a style template) only now, after the
# insert, so it points at the right sheet (it was pushed from index 2 to
# index 3 by the Add() above).
$q3 = $wb.Worksheets.Item(3)              # "2022-Q3" - used as a style template

# Copy the header-row style (bold / centered / bordered) and the index
# column style from the existing "2022-Q3" sheet so formatting matches.
$q3.Range("B1:H1").Copy($newSheet.Range("B1:H1"))
$q3.Range("A2:A2").Copy($newSheet.Range("A2:A26"))

# --- Header row --------------------------------------------------------
$newSheet.Cells.Item(1, 2).Value = "基金代码"
$newSheet.Cells.Item(1, 3).Value = "基金名称"
$newSheet.Cells.Item(1, 4).Value = "基金规模"
$newSheet.Cells.Item(1, 5).Value = "股票总仓位"
$newSheet.Cells.Item(1, 6).Value = "仓位占比"
$newSheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1, 8).Value = "仓位排名"

# --- Data rows -----------------------------------------------------------
# Columns B-G are stored as plain text (matching the source workbook, where
# numeric-looking values such as fund size / position % are text, not
# numbers) - prefix with a single quote to force text entry, then reset the
# number format so no stray "quote prefix" styling is left behind. Column H
# (rank) is a genuine number.
$data = @(
    @("206009", "鹏华新兴产业混合", "43.49", "92.64", "5.36", "2.3311", 3),
    @("001678", "英大国企改革主题股票", "4.20", "92.20", "7.72", "0.3242", 3),
    @("003713", "英大睿盛灵活配置混合A", "2.39", "93.29", "9.84", "0.2352", 1),
    @("003714", "英大睿盛灵活配置混合C", "2.39", "93.29", "9.84", "0.2352", 1),
    @("008132", "鹏华价值驱动混合", "4.41", "93.64", "4.60", "0.2029", 3),
    @("519655", "银河现代服务主题灵活配置混合", "3.36", "87.88", "5.42", "0.1821", 3),
    @("001468", "广发改革先锋灵活配置混合", "5.72", "93.29", "2.47", "0.1413", 8),
    @("011346", "淳厚鑫淳一年持有期混合", "3.34", "78.96", "3.33", "0.1112", 2),
    @("012454", "淳厚鑫悦混合A", "1.82", "85.29", "2.89", "0.0526", 4),
    @("005041", "人保研究精选混合A", "1.09", "87.60", "3.73", "0.0407", 4),
    @("001607", "英大策略优选混合A", "0.59", "93.12", "4.95", "0.0292", 7),
    @("012522", "英大稳固增强核心一年持有混合C", "1.05", "23.17", "2.13", "0.0224", 2),
    @("003446", "英大睿鑫灵活配置混合A", "0.27", "93.18", "8.09", "0.0218", 4),
    @("003447", "英大睿鑫灵活配置混合C", "0.22", "93.18", "8.09", "0.0178", 4),
    @("012455", "淳厚鑫悦混合C", "0.57", "85.29", "2.89", "0.0165", 4),
    @("012521", "英大稳固增强核心一年持有混合A", "0.63", "23.17", "2.13", "0.0134", 2),
    @("519987", "长信恒利优势混合", "0.22", "84.96", "4.59", "0.0101", 8),
    @("006644", "弘毅远方消费升级混合A", "0.39", "84.03", "2.39", "0.0093", 4),
    @("010428", "兴银策略智选混合C", "0.37", "92.14", "2.49", "0.0092", 10),
    @("010427", "兴银策略智选混合A", "0.24", "92.14", "2.49", "0.0060", 10),
    @("001730", "兴银大健康灵活配置混合", "0.15", "91.17", "3.34", "0.0050", 10),
    @("002020", "国都创新驱动灵活配置混合", "0.12", "83.47", "3.02", "0.0036", 9),
    @("005042", "人保研究精选混合C", "0.05", "87.60", "3.73", "0.0019", 4),
    @("014422", "弘毅远方消费升级混合C", "0.06", "84.03", "2.39", "0.0014", 4),
    @("001608", "英大策略优选混合C", "0.02", "93.12", "4.95", "0.0010", 7)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $rec = $data[$i]

    $newSheet.Cells.Item($row, 1).Value = $i

    $newSheet.Cells.Item($row, 2).Value = "'" + $rec[0]
    $newSheet.Cells.Item($row, 2).Style = "Normal"

    $newSheet.Cells.Item($row, 3).Value = "'" + $rec[1]
    $newSheet.Cells.Item($row, 3).Style = "Normal"

    $newSheet.Cells.Item($row, 4).Value = "'" + $rec[2]
    $newSheet.Cells.Item($row, 4).Style = "Normal"

    $newSheet.Cells.Item($row, 5).Value = "'" + $rec[3]
    $newSheet.Cells.Item($row, 5).Style = "Normal"

    $newSheet.Cells.Item($row, 6).Value = "'" + $rec[4]
    $newSheet.Cells.Item($row, 6).Style = "Normal"

    $newSheet.Cells.Item($row, 7).Value = "'" + $rec[5]
    $newSheet.Cells.Item($row, 7).Style = "Normal"

    $newSheet.Cells.Item($row, 8).Value = $rec[6]
}
